# Actualización desde MV -datos-
# Append 10 new daily rows (14-09-2021 .. 27-09-2021) to Sheet1,
# continuing directly after the existing last row (180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("14-09-2021", 37, 7, 9, -2, 0, -2, -2, -40, -13, 7),
    @("15-09-2021", 37, 7, 9, -2, 0, -2, -2, -41, -13, 7),
    @("16-09-2021", 37, 7, 9, -2, 0, -2, -1, -40, -15, 7),
    @("17-09-2021", 37, 7, 9, -2, 0, -2, -1, -40, -15, 7),
    @("20-09-2021", 37, 7, 9, -2, 0, -2, -4, -39, -14, 7),
    @("21-09-2021", 37, 7, 9, -2, 0, -2, -6, -38, -12, 7),
    @("22-09-2021", 37, 7, 9, -2, 0, -2, -4, -37, -14, 7),
    @("23-09-2021", 37, 7, 9, -2, 0, -2, -4, -38, -14, 7),
    @("24-09-2021", 37, 7, 9, -2, 0, -2, -5, -36, -14, 7),
    @("27-09-2021", 37, 7, 9, -2, 0, -2, -5, -35, -16, 7)
)

$startRow = 181

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowData[$col - 1]
    }
}
